$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = 44383
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 16000
$ws.Range("P5").Value = 889

# Row 6
$ws.Range("D6").Value = 44383
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 200

# Row 7
$ws.Range("D7").Value = 44235
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 13000
$ws.Range("P7").Value = 722

# Row 8
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 11000
$ws.Range("L8").Value = 11000
$ws.Range("M8").Value = 11000
$ws.Range("P8").Value = 611

# Row 9
$ws.Range("I9").Value = "Tercera"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = 9000
$ws.Range("P9").Value = 500

# Row 10
$ws.Range("D10").Value = 44396
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 15000
$ws.Range("P10").Value = 833

# Row 11
$ws.Range("D11").Value = 44396
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 150

# Row 12
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 12000
$ws.Range("P12").Value = 667

# Row 13
$ws.Range("D13").Value = 44245
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("P13").Value = 556

# Row 14
$ws.Range("D14").Value = 44229
$ws.Range("I14").Value = "Primera"
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 15000
$ws.Range("P14").Value = 833

# Row 15
$ws.Range("D15").Value = 44249
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 12000
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = 12000
$ws.Range("P15").Value = 667

# Row 16
$ws.Range("D16").Value = 44249
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 10000
$ws.Range("P16").Value = 556

# Row 17
$ws.Range("D17").Value = 44238
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 300

# Row 18
$ws.Range("D18").Value = 44238
$ws.Range("I18").Value = "Segunda"
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = 10000
$ws.Range("P18").Value = 556

# Row 19
$ws.Range("D19").Value = 44238
$ws.Range("I19").Value = "Tercera"
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 8000
$ws.Range("L19").Value = 8000
$ws.Range("M19").Value = 8000
$ws.Range("P19").Value = 444

# Row 20
$ws.Range("D20").Value = 44391
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 15000
$ws.Range("P20").Value = 833
